$d = $word.ActiveDocument

# The document currently ends with an empty paragraph that only holds the
# _GoBack bookmark. We need to insert two brand-new paragraphs of notes
# before it, and add a third note as a new run inside that same paragraph
# (ahead of the bookmark markers).

$lastPara = $d.Paragraphs.Last

$note1 = "Allow starting balances / values for accounts to be set (this could be done with transactions as stated above but it would be easier and more efficient to add a method to be used in the database initializer."
$note2 = "Add sub accounts which can be connected to over arcing accounts such as profit, fixed assets, current assets etc."
$note3 = "Currently the total for all balance sheets can only equal 0 as all transaction have a debt and credit effect of the same value to the accounts which means everything is balanced out, with the addition of profit and loss transactions the profit or loss should be able to be calculated and recorded in a new account i.e. " + [char]0x2018 + "profit" + [char]0x2019 + ", this should be done in addition to the value of one of the effected accounts being higher or lower than the other. This way the balance should be a number other than 0 and this number should equal the value held in the profit account."

# $lastPara.Range stays anchored at the same fixed document position (the
# start of that originally-empty paragraph) across calls, so each
# InsertBefore drops its text immediately before whatever was inserted by
# the previous call. Insert in reverse order (note3, then note2, then
# note1) so the final reading order ends up note1, note2, note3+bookmark.
$lastPara.Range.InsertBefore($note3)
$lastPara.Range.InsertBefore($note2 + "`r")
$lastPara.Range.InsertBefore($note1 + "`r")
